$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.956.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "'2.668.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'566.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "'144.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "'6.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "'3.138.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'26.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.69%  "
$ws.Range("D15").Value = "'60.956.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'2.663.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "'11.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "'351.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'0.528"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.162"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "'8.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("E28").Value = "  +9.26%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +7.26%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'163.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.58%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("E36").Value = "  +6.90%  "
$ws.Range("D37").Value = "'339.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.41%  "
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'4.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.57%  "
$ws.Range("D40").Value = "'0.910"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("D41").Value = "'38.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'5.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("D43").Value = "'0.625"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").Value = "'20.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").Value = "'0.0566"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "'0.0249"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D47").Value = "'133.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'0.0996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'20.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'2.102.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.55%  "
